$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bilateral trade portfolio data (OW-535: align with acuo-data test branch)
$ws.Range("B2").Value = "ACUOSG8745"
$ws.Range("AP2").Value = "p1"

# Collapse the multi-cell selection down to the single active cell (A2),
# matching the refreshed sheet view state.
$ws.Range("A2").Select() | Out-Null
